$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for the Apoe-Lrp1 ligand-receptor pair table
# (values recomputed with the new TPM data; see commit "update scripts wuth new tpm")

# Row 2
$ws.Range("G2").Value = 47.23036199999999
$ws.Range("H2").Value = 141.691086
$ws.Range("I2").Value = 0.3244251370417807
$ws.Range("J2").Value = 0.3244251370417807
$ws.Range("M2").Value = 3.456265333333333
$ws.Range("N2").Value = 10.368796
$ws.Range("O2").Value = 0.009841535807677501
$ws.Range("P2").Value = 0.0098415358076775
$ws.Range("Q2").Value = 163.240662861384
$ws.Range("R2").Value = 1469.165965752456
$ws.Range("S2").Value = 0.003192841603107365
$ws.Range("T2").Value = 0.003192841603107365
# Row 3
$ws.Range("G3").Value = 47.23036199999999
$ws.Range("H3").Value = 141.691086
$ws.Range("I3").Value = 0.3244251370417807
$ws.Range("J3").Value = 0.3244251370417807
$ws.Range("O3").Value = 0.8587907398420774
$ws.Range("P3").Value = 0.8587907398420773
$ws.Range("Q3").Value = 14244.68420078048
$ws.Range("R3").Value = 128202.1578070243
$ws.Range("S3").Value = 0.2786133034634782
$ws.Range("T3").Value = 0.2786133034634782
# Row 4
$ws.Range("G4").Value = 47.23036199999999
$ws.Range("H4").Value = 141.691086
$ws.Range("I4").Value = 0.3244251370417807
$ws.Range("J4").Value = 0.3244251370417807
$ws.Range("O4").Value = 0.1313677243502452
$ws.Range("P4").Value = 0.1313677243502452
$ws.Range("Q4").Value = 2178.984542717046
$ws.Range("R4").Value = 19610.86088445341
$ws.Range("S4").Value = 0.04261899197519517
$ws.Range("T4").Value = 0.04261899197519516
# Row 5
$ws.Range("I5").Value = 0.4188548944674916
$ws.Range("J5").Value = 0.4188548944674916
$ws.Range("M5").Value = 3.456265333333333
$ws.Range("N5").Value = 10.368796
$ws.Range("O5").Value = 0.009841535807677501
$ws.Range("P5").Value = 0.0098415358076775
$ws.Range("Q5").Value = 210.7547868795471
$ws.Range("R5").Value = 1896.793081915924
$ws.Range("S5").Value = 0.0041221754421228
$ws.Range("T5").Value = 0.004122175442122799
# Row 6
$ws.Range("I6").Value = 0.4188548944674916
$ws.Range("J6").Value = 0.4188548944674916
$ws.Range("O6").Value = 0.8587907398420774
$ws.Range("P6").Value = 0.8587907398420773
$ws.Range("S6").Value = 0.3597087047062124
$ws.Range("T6").Value = 0.3597087047062123
# Row 7
$ws.Range("I7").Value = 0.4188548944674916
$ws.Range("J7").Value = 0.4188548944674916
$ws.Range("O7").Value = 0.1313677243502452
$ws.Range("P7").Value = 0.1313677243502452
$ws.Range("S7").Value = 0.05502401431915647
$ws.Range("T7").Value = 0.05502401431915646
# Row 8
$ws.Range("I8").Value = 0.2567199684907278
$ws.Range("J8").Value = 0.2567199684907277
$ws.Range("M8").Value = 3.456265333333333
$ws.Range("N8").Value = 10.368796
$ws.Range("O8").Value = 0.009841535807677501
$ws.Range("P8").Value = 0.0098415358076775
$ws.Range("Q8").Value = 129.1735227680062
$ws.Range("R8").Value = 1162.561704912056
$ws.Range("S8").Value = 0.002526518762447337
$ws.Range("T8").Value = 0.002526518762447336
# Row 9
$ws.Range("I9").Value = 0.2567199684907278
$ws.Range("J9").Value = 0.2567199684907277
$ws.Range("O9").Value = 0.8587907398420774
$ws.Range("P9").Value = 0.8587907398420773
$ws.Range("S9").Value = 0.2204687316723869
$ws.Range("T9").Value = 0.2204687316723868
# Row 10
$ws.Range("I10").Value = 0.2567199684907278
$ws.Range("J10").Value = 0.2567199684907277
$ws.Range("O10").Value = 0.1313677243502452
$ws.Range("P10").Value = 0.1313677243502452
$ws.Range("S10").Value = 0.03372471805589355
$ws.Range("T10").Value = 0.03372471805589354
